$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 205, shifting existing rows 205-208 down to 206-209.
$ws.Rows.Item(205).Insert()

# Populate the newly inserted row 205 with the new record's data.
$ws.Cells.Item(205, 1).Value = 3
$ws.Cells.Item(205, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(205, 3).Value = "Coquimbo"
$ws.Cells.Item(205, 4).Value = 44448
$ws.Cells.Item(205, 5).Value = 5
$ws.Cells.Item(205, 6).Value = 100112003
$ws.Cells.Item(205, 7).Value = "Ajo"
$ws.Cells.Item(205, 8).Value = "Chino"
$ws.Cells.Item(205, 9).Value = "Primera"
$ws.Cells.Item(205, 10).Value = 95
$ws.Cells.Item(205, 11).Value = 15500
$ws.Cells.Item(205, 12).Value = 16000
$ws.Cells.Item(205, 13).Value = 15737
$ws.Cells.Item(205, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(205, 15).Value = "China"
$ws.Cells.Item(205, 16).Value = 1574
$ws.Cells.Item(205, 17).Value = 10
$ws.Cells.Item(205, 18).Value = "Hortaliza"
